$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table only ever had one real row of data left; replace it with a
# placeholder and drop every other company row beneath it.
$ws.Range("A2").Value = "…"
$ws.Range("B2").Value = "…"

# Row 2 becomes the new (and only) data row, so it gets the thick bottom
# border that used to close the table off at row 8.
$ws.Range("A2:B2").Borders.Item(9).LineStyle = 1
$ws.Range("A2:B2").Borders.Item(9).Weight = -4138
$ws.Rows.Item(2).RowHeight = 16

# Clear the now-empty rows 3-8: wipe the values and drop their borders, and
# restore row 8's height/border since it is no longer the closing row.
$ws.Range("A3:B8").ClearContents()
$ws.Range("A3:B8").Borders.LineStyle = -4142
$ws.Rows.Item(8).RowHeight = 15

$ws.Range("B10").Select()
